# Insert a new data row at row 653 (pushing the existing rows 653:667 down
# to 654:668) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(653).Insert()

$ws.Cells.Item(653, 1).Value = 3
$ws.Cells.Item(653, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(653, 3).Value = "Coquimbo"
$ws.Cells.Item(653, 4).Value = 45239
$ws.Cells.Item(653, 5).Value = 5
$ws.Cells.Item(653, 6).Value = 100112043
$ws.Cells.Item(653, 7).Value = "Pepino ensalada"
$ws.Cells.Item(653, 8).Value = "Sin especificar"
$ws.Cells.Item(653, 9).Value = "Primera"
$ws.Cells.Item(653, 10).Value = 65
$ws.Cells.Item(653, 11).Value = 15000
$ws.Cells.Item(653, 12).Value = 15000
$ws.Cells.Item(653, 13).Value = 15000
$ws.Cells.Item(653, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(653, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(653, 16).Value = 250
$ws.Cells.Item(653, 17).Value = 60
$ws.Cells.Item(653, 18).Value = "Hortaliza"
